$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Duration (column C) corrections for existing tasks whose estimated
#    effort changed. Columns B (START DATE) and D (END DATE) are formulas
#    (=previous End / =Start+Duration) so they recompute automatically.
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = 7
$ws.Range("C14").Value = 7
$ws.Range("C15").Value = 7
$ws.Range("C17").Value = 7
$ws.Range("C18").Value = 7
$ws.Range("C20").Value = 16

# ---------------------------------------------------------------------------
# 2. Two new tasks appended to the plan: Edge Node resource handling work.
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "Research efficient way to handle Edge Node resources"
$ws.Range("B23").Formula = "=D22"
$ws.Range("C23").Value = 14
$ws.Range("D23").Formula = "=B23+C23"
$ws.Range("B23").NumberFormat = "d/m;@"
$ws.Range("D23").NumberFormat = "dd/mm/yy;@"

$ws.Range("A24").Value = "Implement Edge Node resource handling"
$ws.Range("B24").Formula = "=D23"
$ws.Range("C24").Value = 14
$ws.Range("D24").Formula = "=B24+C24"
$ws.Range("B24").NumberFormat = "d/m;@"
$ws.Range("D24").NumberFormat = "dd/mm/yy;@"

# New helper column width for the (now more cramped) DURATION column.
$ws.Columns.Item(3).ColumnWidth = 9.7

# ---------------------------------------------------------------------------
# 3. Extend the Gantt bar chart's two series so they cover the new rows.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$series1 = $chart.SeriesCollection().Item(1)
$series1.Formula = "=SERIES('Gantt Chart'!`$B`$1,'Gantt Chart'!`$A`$2:`$A`$24,'Gantt Chart'!`$B`$2:`$B`$24,1)"

$series2 = $chart.SeriesCollection().Item(2)
$series2.Formula = "=SERIES('Gantt Chart'!`$C`$1,'Gantt Chart'!`$A`$2:`$A`$24,'Gantt Chart'!`$C`$2:`$C`$24,2)"

$chart.Refresh()

# ---------------------------------------------------------------------------
# 4. Leave the selection where the edits ended, matching the author's
#    last-touched cell.
# ---------------------------------------------------------------------------
$null = $ws.Range("C24").Select()
